# Fixed variables and query errors in Bread from TC30 to TC47
#
# The CasesTab query (row 2, column B on the "startup" sheet) had a trailing
# `co.cohort_description` / `Cohort` column that is no longer part of the
# query result set. Remove that trailing line (and the stray blank line /
# dangling comma that went with it) so the query ends cleanly after the
# "Response to Treatment" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['West Highland White Terrier'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesTabQuery

# Restore the selection to B2 (it had drifted to B4 / scrolled so row 4 was
# the top-left visible cell).
$ws.Range("B2").Select() | Out-Null

Write-Output "Updated CasesTab query in B2 and reset selection to B2"
